# Add summary-statistic formulas (Max, Min, Mean, Median, Count) to column G
# on every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range('G2').Formula = '=MAX(OFFSET($D$2,0,0,COUNT($D:$D)-1,1))'
    $ws.Range('G3').Formula = '=MIN(OFFSET($D$2,0,0,COUNT($D:$D)-1,1))'
    $ws.Range('G4').Formula = '=IF(G6=0,0,AVERAGE(OFFSET($D$2,0,0,COUNT($D:$D)-1,1)))'
    $ws.Range('G5').Formula = '=IF(G6=0,0,MEDIAN(OFFSET($D$2,0,0,COUNT($D:$D)-1,1)))'
    $ws.Range('G6').Formula = '=COUNT(OFFSET($D$2,0,0,COUNT($D:$D)-1,1))'
}
